# Apply the "make forex gains from dividends tax-free" edit.
#
# Two Foreign-Currency rows (USD lots that originate from dividend payments
# rather than a genuine FOREX purchase, rows 3 and 4) should no longer carry
# a taxable gain: their Gain [EUR] is zeroed out and their Comment is
# replaced with a shared explanatory note. Quantities/gains on the other
# rows are corrected, and the downstream totals (this sheet's summary rows
# plus the mirrored total on the ELSTER summary sheet) are updated to match.

$wb = $excel.ActiveWorkbook

$fx = $wb.Worksheets.Item("Foreign Currencies")

# Row 2: USD lot - quantity correction.
$fx.Range("B2").Value = 1247.91

# Row 3: USD lot received via dividend payment - no longer taxable.
$fx.Range("G3").Value = 0
$fx.Range("H3").Value = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 4: USD lot received via dividend payment - no longer taxable.
$fx.Range("G4").Value = 0
$fx.Range("H4").Value = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 5: USD lot - quantity and recomputed gain.
$fx.Range("B5").Value = 2567.09
$fx.Range("G5").Value = 20.6

# Row 6: USD lot - quantity correction.
$fx.Range("B6").Value = 849.87

# Row 7: USD lot - quantity and recomputed gain.
$fx.Range("B7").Value = 135.13
$fx.Range("G7").Value = -8.83

# Summary rows recomputed after the above changes.
$fx.Range("G9").Value = -43.75
$fx.Range("G10").Value = 20.6
$fx.Range("G11").Value = -64.34999999999999

# ELSTER summary sheet mirrors the FOREX gain/loss total (row 7 = "Zeilen 42
# - 48: Gewinn / Verlust aus Verkauf von Fremdwährungen").
$elster = $wb.Worksheets.Item("ELSTER - Summary")
$elster.Range("C7").Value = -43.75
